$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "diagonal" values for existing rows (one new cell per row)
$ws.Range("J16").Value = 0.2913785036779461
$ws.Range("I17").Value = 0.316329
$ws.Range("H18").Value = 0.3851272760018804
$ws.Range("G19").Value = 0.396329
$ws.Range("F20").Value = 0.5239785795507702
$ws.Range("E21").Value = 0.1852256743764388
$ws.Range("D22").Value = 0.2224853307127136
$ws.Range("C23").Value = 0.1718054801491876

# Add a brand new row 24 - copy A23's formatting (bold/border/centered style)
# down to A24 before setting its value so it picks up the same cell style (s="1").
$ws.Range("A23").Copy($ws.Range("A24"))
$ws.Range("A24").Value = "2025-08-22 00:00:00_diff"
$ws.Range("B24").Value = 0.1722785356205764
